# Update on 12/10/2025 at 8:58pm
#
# The "Table" sheet tracks MMC enrollment/spending by fiscal year. Row 9
# holds the FY26 "YTD" actuals that feed the FY26 full-year estimate in
# row 8 (and cascade up through rows 4-7). This commit refreshes those
# YTD actuals with a newer data pull (7 months reported instead of 3),
# which ripples through all of the dependent formulas automatically.

$wb = $excel.ActiveWorkbook
$wsTable = $wb.Worksheets.Item("Table")
$wsProj  = $wb.Worksheets.Item("projection_data")

# --- Refresh the FY26 (YTD) actuals on the Table sheet -------------------
# B9 = Months Reported, C9 = State Spending YTD, E9 = End of Period Enrollment YTD
$wsTable.Range("B9").Value = 7
$wsTable.Range("C9").Value = 31158106
$wsTable.Range("E9").Value = 4419503

# D9 (Average Annual Enrollment, YTD) isn't a live formula - it mirrors the
# freshly recalculated "per-month" run-rate in E8, same as the prior pull.
$wsTable.Range("D9").Value = $wsTable.Range("E8").Value2

# Column B ("Months Reported") was hidden; unhide it so the new figure is visible.
$wsTable.Columns("B").Hidden = $false

# --- Restore the user's on-screen selections ------------------------------
# Update the selection on projection_data first (without leaving it the
# active tab), then reselect the Table sheet and its new active cell.
[void]$wsProj.Range("F23").Select()
[void]$wsTable.Select()
[void]$wsTable.Range("G30").Select()
